$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Data Validation" column before the existing "Data Testing" column (D),
# shifting the old D column (Data Testing) to E.
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column D
$ws.Range("D1").Value = "Data Validation"

# Update row 2 values
$ws.Range("B2").Value = 40320
$ws.Range("C2").Value = 22579
$ws.Range("D2").Value = 5645
$ws.Range("E2").Value = 12096

# Update row 3 values
$ws.Range("B3").Value = 40320
$ws.Range("C3").Value = 22579
$ws.Range("D3").Value = 5645
$ws.Range("E3").Value = 12096
